$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the NATMI ligand-receptor edge table (rows 2-10) with the
# recomputed values ("Natmi following Dr Hou advice"): the number of
# ligand-/receptor-expressing cells changed from 1 to 3 for every row,
# which ripples through the total-expression, specificity and edge-weight
# columns. Columns A-D, F and L are unaffected.

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 6.138059000000001
$ws.Range("H2").Value = 18.414177
$ws.Range("I2").Value = 0.08535364925338249
$ws.Range("J2").Value = 0.08535364925338247
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 112.513392
$ws.Range("N2").Value = 337.540176
$ws.Range("O2").Value = 0.3275312977368564
$ws.Range("P2").Value = 0.3275312977368564
$ws.Range("Q2").Value = 690.6138383861281
$ws.Range("R2").Value = 6215.524545475152
$ws.Range("S2").Value = 0.02795599150653683
$ws.Range("T2").Value = 0.02795599150653682

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 6.138059000000001
$ws.Range("H3").Value = 18.414177
$ws.Range("I3").Value = 0.08535364925338249
$ws.Range("J3").Value = 0.08535364925338247
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 106.314466
$ws.Range("N3").Value = 318.943398
$ws.Range("O3").Value = 0.3094859589441663
$ws.Range("P3").Value = 0.3094859589441664
$ws.Range("Q3").Value = 652.5644648614941
$ws.Range("R3").Value = 5873.080183753446
$ws.Range("S3").Value = 0.0264157559885671
$ws.Range("T3").Value = 0.0264157559885671

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 6.138059000000001
$ws.Range("H4").Value = 18.414177
$ws.Range("I4").Value = 0.08535364925338249
$ws.Range("J4").Value = 0.08535364925338247
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 124.6916553333333
$ws.Range("N4").Value = 374.074966
$ws.Range("O4").Value = 0.3629827433189773
$ws.Range("P4").Value = 0.3629827433189773
$ws.Range("Q4").Value = 765.3647372436649
$ws.Range("R4").Value = 6888.282635192983
$ws.Range("S4").Value = 0.03098190175827855
$ws.Range("T4").Value = 0.03098190175827855

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 37.15353
$ws.Range("H5").Value = 111.46059
$ws.Range("I5").Value = 0.5166436764692264
$ws.Range("J5").Value = 0.5166436764692264
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 112.513392
$ws.Range("N5").Value = 337.540176
$ws.Range("O5").Value = 0.3275312977368564
$ws.Range("P5").Value = 0.3275312977368564
$ws.Range("Q5").Value = 4180.269685073759
$ws.Range("R5").Value = 37622.42716566384
$ws.Range("S5").Value = 0.1692169738215063
$ws.Range("T5").Value = 0.1692169738215063

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 37.15353
$ws.Range("H6").Value = 111.46059
$ws.Range("I6").Value = 0.5166436764692264
$ws.Range("J6").Value = 0.5166436764692264
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 106.314466
$ws.Range("N6").Value = 318.943398
$ws.Range("O6").Value = 0.3094859589441663
$ws.Range("P6").Value = 0.3094859589441664
$ws.Range("Q6").Value = 3949.95770196498
$ws.Range("R6").Value = 35549.61931768482
$ws.Range("S6").Value = 0.1598939636445181
$ws.Range("T6").Value = 0.1598939636445182

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 37.15353
$ws.Range("H7").Value = 111.46059
$ws.Range("I7").Value = 0.5166436764692264
$ws.Range("J7").Value = 0.5166436764692264
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 124.6916553333333
$ws.Range("N7").Value = 374.074966
$ws.Range("O7").Value = 0.3629827433189773
$ws.Range("P7").Value = 0.3629827433189773
$ws.Range("Q7").Value = 4632.73515717666
$ws.Range("R7").Value = 41694.61641458994
$ws.Range("S7").Value = 0.1875327390032019
$ws.Range("T7").Value = 0.187532739003202

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 28.621669
$ws.Range("H8").Value = 85.86500699999999
$ws.Range("I8").Value = 0.3980026742773913
$ws.Range("J8").Value = 0.3980026742773912
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 112.513392
$ws.Range("N8").Value = 337.540176
$ws.Range("O8").Value = 0.3275312977368564
$ws.Range("P8").Value = 0.3275312977368564
$ws.Range("Q8").Value = 3220.321063891247
$ws.Range("R8").Value = 28982.88957502123
$ws.Range("S8").Value = 0.1303583324088133
$ws.Range("T8").Value = 0.1303583324088133

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 28.621669
$ws.Range("H9").Value = 85.86500699999999
$ws.Range("I9").Value = 0.3980026742773913
$ws.Range("J9").Value = 0.3980026742773912
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 106.314466
$ws.Range("N9").Value = 318.943398
$ws.Range("O9").Value = 0.3094859589441663
$ws.Range("P9").Value = 0.3094859589441664
$ws.Range("Q9").Value = 3042.897455763753
$ws.Range("R9").Value = 27386.07710187378
$ws.Range("S9").Value = 0.1231762393110811
$ws.Range("T9").Value = 0.1231762393110811

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 28.621669
$ws.Range("H10").Value = 85.86500699999999
$ws.Range("I10").Value = 0.3980026742773913
$ws.Range("J10").Value = 0.3980026742773912
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 124.6916553333333
$ws.Range("N10").Value = 374.074966
$ws.Range("O10").Value = 0.3629827433189773
$ws.Range("P10").Value = 0.3629827433189773
$ws.Range("Q10").Value = 3568.883286012751
$ws.Range("R10").Value = 32119.94957411476
$ws.Range("S10").Value = 0.1444681025574968
$ws.Range("T10").Value = 0.1444681025574968
